# Auto-upload VRF Excel file
# Add a new worksheet named "new" at the end of the workbook with the
# standard VRF header row (Outdoor/Indoor Model, Quantity, Serial(s)).

$wb = $excel.ActiveWorkbook

# Add a brand-new worksheet and move it to become the very last tab.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "new"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)

# Re-resolve the worksheet by name: after Move() shuffles tab order, any
# previously captured reference can point at whatever now sits at the old
# index, so look the sheet up fresh by its (now final) name.
$ws = $wb.Worksheets.Item("new")

# Header row values first ...
$headers = @("Outdoor Model", "Outdoor Quantity", "Outdoor Serial(s)", "Indoor Model", "Indoor Quantity", "Indoor Serial(s)")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ... then formatting in one shot per property over the whole header range,
# so every header cell lands on the same style index (bold, thin box
# border, centered / top-aligned) instead of Excel minting a fresh style
# per incremental per-cell tweak.
$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1        # xlContinuous
$headerRange.Borders.Weight = 2           # xlThin
